# Weekly price-sheet refresh: insert this week's Cilantro price observation
# at the top of the historical price list (row 83, just under the header
# block) and push the previous observations down by one row.
#
# The workbook lists individual price observations for
# "Terminal La Palmera de La Serena" Cilantro sorted with the newest
# observation first starting at row 83. Adding this week's observation
# means inserting a new row 83 and shifting every existing row (83..106)
# down by one (83->84, ..., 106->107); the bottom-most historical row
# (old row 106) ends up at row 107.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 83; Excel shifts rows 83:106 down to 84:107
# and the new row inherits the number formatting (e.g. the date style on
# column D) from the row that used to be there.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row with this week's observation.
$ws.Cells.Item(83, 1).Value = 8
$ws.Cells.Item(83, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(83, 3).Value = "Coquimbo"
$ws.Cells.Item(83, 4).Value = 44551
$ws.Cells.Item(83, 5).Value = 4
$ws.Cells.Item(83, 6).Value = 100112040
$ws.Cells.Item(83, 7).Value = "Cilantro"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 3200
$ws.Cells.Item(83, 11).Value = 2000
$ws.Cells.Item(83, 12).Value = 2500
$ws.Cells.Item(83, 13).Value = 2250
$ws.Cells.Item(83, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(83, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(83, 16).Value = 1500
$ws.Cells.Item(83, 17).Value = 1.5
$ws.Cells.Item(83, 18).Value = "Hortaliza"
